$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the row stays text-formatted (values look numeric but must remain strings)
$ws.Range("A2:K2").NumberFormat = "@"

# Overwrite row 2 with what used to be row 3's data (the " Oct 7 2020" match vs Chennai Super Kings)
$ws.Range("A2").Value = " Oct 7 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "KKR won by 10 runs"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Eoin Morgan "
$ws.Range("G2").Value = "7"
$ws.Range("H2").Value = "10"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "70.00"

# Remove the old rows 3 and 4 entirely, shrinking the used range down to A1:K2
$ws.Rows("3:4").Delete()
